$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    ,@('Sep 09 2017', 'negative', -0.1, 'https://twitter.com/statuses/905906972393963520', '#0daytoday #Apache Struts 2 REST Plugin XStream Remote Code Execution #Exploit #0day https://t.co/rwaxBrUvuC')
    ,@('Sep 09 2017', 'negative', -0.1, 'https://twitter.com/statuses/905906732098088966', '#0daytoday #SourceTree Remote Code Execution #Exploit https://t.co/cGo5U8zYGm')
    ,@('Sep 09 2017', 'negative', -0.1, 'https://twitter.com/statuses/905878605128884225', '#0daytoday #Gh0st Client - Buffer Overflow Exploit [remote #exploits #0day #Exploit] https://t.co/ce22CcFpwa')
    ,@('Sep 09 2017', 'negative', -0.1, 'https://twitter.com/statuses/905878602645897217', '#0daytoday #PlugX Controller Stack Overflow Exploit [remote #exploits #0day #Exploit] https://t.co/jcciGj0kBX')
    ,@('Sep 09 2017', 'negative', -0.1, 'https://twitter.com/statuses/905750181722947584', '#0daytoday #Apache #Struts 2.5 - Remote Code Execution #0day #Exploit https://t.co/NRIvcrWlEq')
    ,@('Sep 09 2017', 'negative', -0.05, 'https://twitter.com/statuses/905728881197408257', '#0daytoday #WordPress Gym Management System Code Execution / Remote Cross Site Scripting Vulnerabil [#0day #Exploit] https://t.co/6BbjaqVOyR')
    ,@('Sep 09 2017', 'negative', -0.1, 'https://twitter.com/statuses/905513948224438272', '#0daytoday #Tor - Linux Sandbox Breakout via X11 Exploit [remote #exploits #0day #Exploit] https://t.co/8bCWbuT4Uv')
    ,@('Sep 09 2017', 'negative', -0.1, 'https://twitter.com/statuses/905505647285813252', '#0daytoday #Jungo DriverWizard WinDriver - Kernel Pool Overflow Exploit [remote #exploits #0day #Exploit] https://t.co/jc9QMlrawK')
    ,@('Sep 09 2017', 'negative', -0.1, 'https://twitter.com/statuses/905355158887907328', '#0daytoday #Samsung Internet Browser - SOP Bypass Exploit [remote #exploits #0day #Exploit] https://t.co/ka6tQuu1ZI')
    ,@('Sep 09 2017', 'negative', -0.1, 'https://twitter.com/statuses/905079577915088897', '#0daytoday #NEC EXPRESS CLUSTER clpwebmc Remote #Root #Exploit #0day https://t.co/4i3pB8Hwp8')
    ,@('Sep 09 2017', 'negative', -0.1, 'https://twitter.com/statuses/904803563603406850', '#0daytoday #RubyGems < 2.6.13 - Arbitrary File Overwrite #Exploit https://t.co/sM5XxS3Rxt')
    ,@('Sep 09 2017', 'negative', -0.1, 'https://twitter.com/statuses/904743627418513408', '#0daytoday #Wireless Repeater BE126 - Remote Code Execution Exploit [webapps #exploits #0day #Exploit] https://t.co/MJBKKmDAxw')
    ,@('Aug 08 2017', 'negative', -0.1, 'https://twitter.com/statuses/902990838938058752', '#0daytoday #Joomla Joomanager 2.0.0 Component - Arbitrary File Download Vulnerability [#0day #Exploit] https://t.co/OJCfpL4p9I')
)

$r = 4
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]

    $ws.Hyperlinks.Add($ws.Cells.Item($r, 4), $row[3]) | Out-Null
    $ws.Cells.Item($r, 4).Font.Underline = $true
    $ws.Cells.Item($r, 4).Font.Color = 16711680

    $r = $r + 1
}
